# chore: modify example tables for test
#
# Adds a new "random" (int[]) field to the "struct Reward" example table on
# 工作表1, and updates the active sheet / cell-selection bookkeeping so that
# 工作表1 (instead of 工作表2) is the active tab, matching the new selections
# recorded in each sheet's view.

$wb = $excel.ActiveWorkbook

# --- 工作表1 ("struct Reward" example) ---------------------------------
$ws1 = $wb.Worksheets.Item("工作表1")

# New row entries for the extra "random" field of type "int[]"
$ws1.Range("D3").Value = "random"
$ws1.Range("D4").Value = "int[]"

# --- 工作表2 --------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("工作表2")
$ws2.Select()
$ws2.Range("A10").Select()

# --- Make 工作表1 the active sheet, with its new selection ---------------
$ws1.Select()
$ws1.Range("I19").Select()
